$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the operating-system skills (Windows, Linux, Mac OS) from the
# "DevOps" section …
$ws.Rows("54:56").Delete()

# … and give them their own "Betriebssystem" section, placed right after
# the "AWS" block and before "Standardsoftware".
$ws.Rows("83:86").Insert()

$ws.Range("A83").Value = "Betriebssystem"
$ws.Range("D83").Value = 0

$ws.Range("B84").Value = "Windows"
$ws.Range("C84").Value = 5

$ws.Range("B85").Value = "Linux"
$ws.Range("C85").Value = 4

$ws.Range("B86").Value = "Mac OS"
$ws.Range("C86").Value = 3

# Update the "Aktualisiert" date.
$ws.Range("B2").Value = "03.11.2023"

# Match the author's final cursor position.
$ws.Range("I11").Select()
